$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '61.370.15'
$ws.Cells.Item(2, 5).Value = '  +0.08%  '

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.439.88'
$ws.Cells.Item(3, 5).Value = '  +0.21%  '

$ws.Cells.Item(4, 5).Value = '  +0.07%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '575.08'
$ws.Cells.Item(5, 5).Value = '  +0.16%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '144.45'
$ws.Cells.Item(6, 5).Value = '  +3.43%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '3.440.71'
$ws.Cells.Item(7, 5).Value = '  +0.27%  '

$ws.Cells.Item(8, 5).Value = '  +0.07%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.478'
$ws.Cells.Item(9, 5).Value = '  +1.12%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '7.61'
$ws.Cells.Item(10, 5).Value = '  -1.17%  '

$ws.Cells.Item(11, 5).Value = '  +3.16%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.390'
$ws.Cells.Item(12, 5).Value = '  +1.58%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '4.030.21'
$ws.Cells.Item(13, 5).Value = '  +0.45%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '28.73'
$ws.Cells.Item(14, 5).Value = '  +7.37%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.123'
$ws.Cells.Item(15, 5).Value = '  -0.64%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '0.0000173'
$ws.Cells.Item(16, 5).Value = '  +0.22%  '

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '3.451.11'
$ws.Cells.Item(17, 5).Value = '  +0.85%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '61.482.96'
$ws.Cells.Item(18, 5).Value = '  +0.15%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '6.37'
$ws.Cells.Item(19, 5).Value = '  +6.98%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '14.46'
$ws.Cells.Item(20, 5).Value = '  +3.60%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '9.36'
$ws.Cells.Item(21, 5).Value = '  -1.56%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '401.26'
$ws.Cells.Item(22, 5).Value = '  +5.07%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '0.565'
$ws.Cells.Item(23, 5).Value = '  +1.29%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '73.86'
$ws.Cells.Item(24, 5).Value = '  +3.14%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '1.00'
$ws.Cells.Item(25, 5).Value = '  +0.41%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '0.0000122'
$ws.Cells.Item(26, 5).Value = '  -0.93%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '3.592.26'
$ws.Cells.Item(27, 5).Value = '  +0.95%  '

$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '0.179'
$ws.Cells.Item(28, 5).Value = '  +0.71%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '7.58'
$ws.Cells.Item(29, 5).Value = '  -1.01%  '

$ws.Cells.Item(30, 5).Value = '  +0.07%  '

$ws.Cells.Item(31, 5).Value = '  -5.47%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '8.24'
$ws.Cells.Item(32, 5).Value = '  +0.93%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '2.18'
$ws.Cells.Item(33, 5).Value = '  +1.27%  '

$ws.Cells.Item(34, 5).Value = '  -0.11%  '

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '23.95'
$ws.Cells.Item(35, 5).Value = '  +0.59%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '7.04'
$ws.Cells.Item(36, 5).Value = '  +1.04%  '

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '3.469.25'
$ws.Cells.Item(37, 5).Value = '  +0.62%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '5.14'
$ws.Cells.Item(38, 5).Value = '  -2.00%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '1.55'
$ws.Cells.Item(39, 5).Value = '  -1.21%  '

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '166.81'
$ws.Cells.Item(40, 5).Value = '  +0.02%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.0792'
$ws.Cells.Item(41, 5).Value = '  +1.27%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '27.42'
$ws.Cells.Item(42, 5).Value = '  +2.96%  '

$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.804'
$ws.Cells.Item(43, 5).Value = '  +2.37%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '4.52'
$ws.Cells.Item(44, 5).Value = '  +2.20%  '

$ws.Cells.Item(45, 5).Value = '  +0.17%  '

$ws.Cells.Item(46, 2).Value = 'OKB'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '42.37'
$ws.Cells.Item(46, 5).Value = '  +0.78%  '

$ws.Cells.Item(47, 2).Value = 'Stacks'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '1.73'
$ws.Cells.Item(47, 5).Value = '  -0.52%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '2.616.16'
$ws.Cells.Item(48, 5).Value = '  -0.84%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.15'
$ws.Cells.Item(49, 5).Value = '  -3.46%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '6.97'
$ws.Cells.Item(50, 5).Value = '  +2.29%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '23.23'
$ws.Cells.Item(51, 5).Value = '  -3.03%  '
